$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing F11 text (shared string index 43) ---
$f11text = "日语近义词/提供20个和主题内容相似的日语单词，提供例句和中文翻译，讲解语法，具体说明使用上的差别。`n回答样式:[XX][XX]...(if more than 1 items in each keyword related content)...Final output are in the following format:     - 段落 1     - 段落 2     - 段落 3`n"
$ws.Range("F11").Value = $f11text

# --- Add two new rows of content below row 11 ---
$f12text = "人物基本情况介绍`n回答样式:[XX][XX]...(if more than 1 items in each keyword related content)...Final output are in the following format:     - 段落 1     - 段落 2     - 段落 3`n"
$f13text = "国家基本情况介绍`n回答样式:[XX][XX]...(In well-structured way. If more than 1 items then use more [] to input content for each keyword related content)...Final output are in the following format:     - 段落 1     - 段落 2     - 段落 3`n"

$ws.Range("F12").Value = $f12text
$ws.Range("F13").Value = $f13text

# Match the wrap/vertical-center formatting already used by F11 on the new cells
$ws.Range("F12").WrapText = $true
$ws.Range("F12").VerticalAlignment = -4108
$ws.Range("F13").WrapText = $true
$ws.Range("F13").VerticalAlignment = -4108

# --- Row heights ---
$ws.Rows.Item(11).RowHeight = 129.6
$ws.Rows.Item(12).RowHeight = 100.8
$ws.Rows.Item(13).RowHeight = 115.2

# Row 10's content no longer wraps, so let it shrink back to the standard height
$ws.Range("F10").WrapText = $false
$ws.Rows.Item(10).AutoFit() | Out-Null

# --- Column widths: widen column F to fit the long text, keep others standard ---
$ws.Columns.Item(6).ColumnWidth = 49

# --- Update selection to follow the new last cell, matching the authored file ---
$ws.Range("F14").Select() | Out-Null
